$wb = $excel.ActiveWorkbook

# "Michael" is the 3rd tab (Arpit, Brody, Michael, Sakshyam, Vasilis, Yong)
$wsMichael = $wb.Worksheets.Item("Michael")

# Update Michael's weekly report: "Issues encountered:" (B9) and
# "Issues resolved:" (B10) with this sprint's GitKraken/GitHub story.
$wsMichael.Range("B9").Value = "GitKracken was not working; therefore, it was not allowing me to upload any of my documents to github or pull from the github!"
$wsMichael.Range("B10").Value = "Vasilis and Brody helped me out in learning how to pull and push to github using the Command Prompt. "

# Switch focus to the Michael sheet (becomes the selected/active tab)
# and leave the cursor on A11, matching the saved selection state.
$wsMichael.Activate()
$wsMichael.Range("A11").Select()
